$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the three now-completed features with their completed version and hide their rows.
$ws.Range("C13").Value = "1.5.6"
$ws.Range("C14").Value = "1.5.6"
$ws.Range("C15").Value = "1.5.6"

$ws.Rows.Item(13).Hidden = $true
$ws.Rows.Item(14).Hidden = $true
$ws.Rows.Item(15).Hidden = $true

# Add new feature requests.
$ws.Range("A16").Value = "Horizontal custom rolls"
$ws.Range("B16").Value = "Having a bunch of the same icon with little text is not an efficient use of space. Convert to horizontal and remove icon."
$ws.Range("D16").Value = "Weston Fiala"

$ws.Range("A17").Value = "Add Genesys dice"
$ws.Range("B17").Value = "The game Genesys uses dice. Add them in."
$ws.Range("D17").Value = "Joseph Thompson - Store Review"

$ws.Range("A18").Value = "Add Fantasy Flight Star Wars dice"
$ws.Range("B18").Value = "The game Fantasy Flight Star Wars uses dice. Add them in."
$ws.Range("D18").Value = "Joseph Thompson - Store Review"

# Update the selection to match the author's final cursor position.
$ws.Range("B21").Select()
